$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data (rows 2-5) is re-sorted by Fecha (column D) ascending, as
# part of the weekly fruit/vegetable price refresh. Only the numeric columns
# (D, M, N, O, P, S) actually move between rows; the descriptive/text columns
# are identical across all four rows, so we just rewrite the numeric values in
# their new, date-sorted order without touching the string cells.

$ws.Range("D2").Value = 44250
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 14500
$ws.Range("S2").Value = 806

$ws.Range("D3").Value = 44252
$ws.Range("M3").Value = 120
$ws.Range("N3").Value = 13000
$ws.Range("O3").Value = 14000
$ws.Range("P3").Value = 13500
$ws.Range("S3").Value = 750

$ws.Range("D4").Value = 44253
$ws.Range("M4").Value = 160
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 14500
$ws.Range("S4").Value = 806

$ws.Range("D5").Value = 44257
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 14500
$ws.Range("S5").Value = 806
